$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 187.6589263439178
$ws.Range("C2").Value = 1.260957667245832
$ws.Range("D2").Value = 1.597512054443359
$ws.Range("E2").Value = 0.05880919239974763
$ws.Range("B3").Value = 371.0399147510528
$ws.Range("C3").Value = 2.247013835908929
$ws.Range("D3").Value = 2.100662755966186
$ws.Range("E3").Value = 0.08295081311260794
$ws.Range("B4").Value = 743.5080466270447
$ws.Range("C4").Value = 2.095083124156394
$ws.Range("D4").Value = 1.839073610305786
$ws.Range("E4").Value = 0.2511388520340822
$ws.Range("B5").Value = 192.3955340385437
$ws.Range("C5").Value = 0.8502665275713333
$ws.Range("D5").Value = 1.85203332901001
$ws.Range("E5").Value = 0.2128062639015629
$ws.Range("B6").Value = 374.0436499118805
$ws.Range("C6").Value = 1.454906537405201
$ws.Range("D6").Value = 1.804590368270874
$ws.Range("E6").Value = 0.2036292224368807
$ws.Range("B7").Value = 738.6901639461518
$ws.Range("C7").Value = 2.566665855546967
$ws.Range("D7").Value = 1.748728656768799
$ws.Range("E7").Value = 0.1248934444475369
$ws.Range("B8").Value = 190.1935725212097
$ws.Range("C8").Value = 0.9022416757328573
$ws.Range("D8").Value = 1.666329383850098
$ws.Range("E8").Value = 0.05473784477056894
$ws.Range("B9").Value = 371.3957646846771
$ws.Range("C9").Value = 2.046691190822151
$ws.Range("D9").Value = 1.677208185195923
$ws.Range("E9").Value = 0.06369654532750987
$ws.Range("B10").Value = 738.8924477577209
$ws.Range("C10").Value = 3.163684878407548
$ws.Range("D10").Value = 1.656140327453613
$ws.Range("E10").Value = 0.03695500140200868
$ws.Range("B11").Value = 314.7178371429443
$ws.Range("C11").Value = 2.13554558872918
$ws.Range("D11").Value = 1.836596202850342
$ws.Range("E11").Value = 0.2289297721418531
$ws.Range("B12").Value = 615.3876048088074
$ws.Range("C12").Value = 2.32133518620277
$ws.Range("D12").Value = 1.804598760604858
$ws.Range("E12").Value = 0.2103760257686744
$ws.Range("B13").Value = 1215.946527576446
$ws.Range("C13").Value = 5.908470233480413
$ws.Range("D13").Value = 1.690570545196533
$ws.Range("E13").Value = 0.0655149278540159
$ws.Range("B14").Value = 311.0300178050995
$ws.Range("C14").Value = 2.47344466598569
$ws.Range("D14").Value = 1.716483974456787
$ws.Range("E14").Value = 0.1534862982390311
$ws.Range("B15").Value = 615.7575721740723
$ws.Range("C15").Value = 3.853422812225658
$ws.Range("D15").Value = 1.690686511993408
$ws.Range("E15").Value = 0.146666769289683
$ws.Range("B16").Value = 1219.672704648972
$ws.Range("C16").Value = 3.93503618757379
$ws.Range("D16").Value = 1.669526672363281
$ws.Range("E16").Value = 0.02889518210339775
$ws.Range("B17").Value = 313.5077887058258
$ws.Range("C17").Value = 4.959640734812297
$ws.Range("D17").Value = 1.818214702606201
$ws.Range("E17").Value = 0.2103209486118717
$ws.Range("B18").Value = 608.0419114112854
$ws.Range("C18").Value = 8.698357761212964
$ws.Range("D18").Value = 1.676743793487549
$ws.Range("E18").Value = 0.1597758415913075
$ws.Range("B19").Value = 982.3444380283356
$ws.Range("C19").Value = 28.61750343890318
$ws.Range("D19").Value = 1.125482702255249
$ws.Range("E19").Value = 0.1689588076926846
